$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new text value into I10 (creates a new shared string entry "fadfdsfas")
$ws.Range("I10").Value = "fadfdsfas"

# Update the selection to match the target state (N10)
$ws.Range("N10").Select()
